$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Report Generated On" timestamp
$ws.Range("D5").Value = "Report Generated On: 08/26/2025 09:59 AM"

# Update Total Billed Amount
$ws.Range("C8").Value = 4918.07

# Clear Scope ID # value
$ws.Range("G10").Value = ""

# Update line item pricing values
$ws.Range("H16").Value = 478.55
$ws.Range("H17").Value = 238.2
$ws.Range("H18").Value = 1905.6
$ws.Range("H19").Value = 476.4
$ws.Range("H20").Value = 238.2
$ws.Range("H21").Value = 1581.12

# Update TOTAL
$ws.Range("H22").Value = 4918.07
